# Update "想去人数" (number of people who want to go) values in column F
# for sheets "展览" and "全部类型" (rows 5,6,8,9,10,11,12,13).

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 523
    6  = 6970
    8  = 158
    9  = 1050
    10 = 415
    11 = 144
    12 = 188
    13 = 595
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
